$d = $word.ActiveDocument

function Insert-LineBreak($pos) {
    # Inserts a manual line break (<w:br/>) as its own run at $pos,
    # returns the position immediately after the break.
    $r = $d.Range($pos, $pos)
    $r.InsertBreak(6) | Out-Null
    return $pos + 1
}

function Insert-StyledRun($pos, $text, $style) {
    # Inserts $text at $pos as its own run, then stamps the character
    # style on that freshly inserted run (forces a run split instead of
    # merging into a neighbouring run). Returns the position right after
    # the inserted text.
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text) | Out-Null
    $r.Style = $style
    return $r.End
}

function Find-Start($searchText, $fromPos) {
    # Finds $searchText starting at $fromPos, returns the Range.
    $r = $d.Range($fromPos, $d.Content.End)
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find: $searchText"
    }
    return $r
}

# --- def p1(): -> add a one-line docstring right after the signature ---
function Add-OneLineDocstring($signature, $docText) {
    $r = Find-Start $signature 0
    $pos = $r.End
    $pos = Insert-LineBreak $pos
    $pos = Insert-StyledRun $pos "    " "NormalTok"
    $pos = Insert-StyledRun $pos $docText "CommentTok"
}

# --- def _p2_p3_res(): -> add the multi-line docstring right after the
# signature, reusing the break that is already there and adding one new
# break right before the (preserved) original body text. ---
function Add-MultiLineDocstring($fromPos) {
    $r = Find-Start "_p2_p3_res():" $fromPos
    $pos = $r.End + 1

    $pos = Insert-StyledRun $pos "    " "NormalTok"
    $pos = Insert-StyledRun $pos '"""' "CommentTok"
    $pos = Insert-LineBreak $pos
    $pos = Insert-StyledRun $pos "    Collapse certain branches and calculate resistances for problems 2 and 3." "CommentTok"
    $pos = Insert-LineBreak $pos
    $pos = Insert-StyledRun $pos "    Reduces repeat calculations." "CommentTok"
    $pos = Insert-LineBreak $pos
    $pos = Insert-LineBreak $pos
    $pos = Insert-StyledRun $pos "    :return tuple: collapsed resistances" "CommentTok"
    $pos = Insert-LineBreak $pos
    $pos = Insert-StyledRun $pos '    """' "CommentTok"
    $pos = Insert-LineBreak $pos

    return $r.End
}

Add-OneLineDocstring "p1():" '"""Problem 1"""'
$afterFirst = Add-MultiLineDocstring 0
Add-MultiLineDocstring $afterFirst
Add-OneLineDocstring "p2():" '"""Problem 2"""'
Add-OneLineDocstring "p3():" '"""Problem 3"""'
Add-OneLineDocstring "p4():" '"""Problem 4"""'
